$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows 2..51: set IPC PO (C) to 0, DELTA (D) to -B, DELTA^2 (E) to B^2
for ($r = 2; $r -le 51; $r++) {
    $b = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r, 3).Value = 0
    $ws.Cells.Item($r, 4).Value = -$b
    $ws.Cells.Item($r, 5).Value = $b * $b
}

# Row 52: TOTAL row -> C52 = sum of D2:D51, E52 = sum of E2:E51
$ws.Cells.Item(52, 3).Value = $ws.Application.WorksheetFunction.Sum($ws.Range("D2:D51"))
$ws.Cells.Item(52, 5).Value = $ws.Application.WorksheetFunction.Sum($ws.Range("E2:E51"))

# Row 53: MSE row -> E53 = average of E2:E51
$ws.Cells.Item(53, 5).Value = $ws.Application.WorksheetFunction.Average($ws.Range("E2:E51"))
